$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scheduler app appended new login-timestamp records since the
# workbook was last saved; the most recent timestamp for this student's
# session (row 2) is now 03/30/2020 23:23:53.
$ws.Range("B2").Value2 = "03/30/2020 23:23:53"
